$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("B1").Value = 45678
$ws.Range("C1").Value = 45685
$ws.Range("D1").Value = 45692
$ws.Range("E1").Value = 45699
$ws.Range("F1").Value = 45706
$ws.Range("G1").Value = 45713
$ws.Range("H1").Value = 45720
$ws.Range("I1").Value = 45727
$ws.Range("J1").Value = 45734
$ws.Range("K1").Value = 45741
$ws.Range("L1").Value = 45748

# Row 2
$ws.Range("B2").Value = 44.1
$ws.Range("C2").Value = 43
$ws.Range("D2").Value = 42.5
$ws.Range("E2").Value = 42.4
$ws.Range("F2").Value = 42.7
$ws.Range("G2").Value = 43.5
$ws.Range("H2").Value = 44.7
$ws.Range("I2").Value = 46.3
$ws.Range("J2").Value = 48.2
$ws.Range("K2").Value = 50.4
$ws.Range("L2").Value = 52.9

# Row 3
$ws.Range("B3").Value = 43.2
$ws.Range("C3").Value = 42
$ws.Range("D3").Value = 41.1
$ws.Range("E3").Value = 40.6
$ws.Range("F3").Value = 40.4
$ws.Range("G3").Value = 40.7
$ws.Range("H3").Value = 41.3
$ws.Range("I3").Value = 42.3
$ws.Range("J3").Value = 43.8
$ws.Range("K3").Value = 45.5
$ws.Range("L3").Value = 47.6

# Row 4
$ws.Range("B4").Value = 41.9
$ws.Range("C4").Value = 39.2
$ws.Range("D4").Value = 37.1
$ws.Range("E4").Value = 35.6
$ws.Range("F4").Value = 35
$ws.Range("G4").Value = 35.3
$ws.Range("H4").Value = 36.4
$ws.Range("I4").Value = 38.3
$ws.Range("J4").Value = 40.9
$ws.Range("K4").Value = 44.1
$ws.Range("L4").Value = 47.5

# Row 5
$ws.Range("B5").Value = 66.40000000000001
$ws.Range("C5").Value = 66.59999999999999
$ws.Range("D5").Value = 66.40000000000001
$ws.Range("E5").Value = 65.7
$ws.Range("F5").Value = 64.8
$ws.Range("G5").Value = 63.5
$ws.Range("H5").Value = 62.2
$ws.Range("I5").Value = 60.8
$ws.Range("J5").Value = 59.6
$ws.Range("K5").Value = 58.6
$ws.Range("L5").Value = 58

# Row 6
$ws.Range("B6").Value = 43.3
$ws.Range("C6").Value = 42.6
$ws.Range("D6").Value = 42.5
$ws.Range("E6").Value = 43
$ws.Range("F6").Value = 43.8
$ws.Range("G6").Value = 45
$ws.Range("H6").Value = 46.3
$ws.Range("I6").Value = 47.7
$ws.Range("J6").Value = 49.3
$ws.Range("K6").Value = 50.9
$ws.Range("L6").Value = 52.5

# Row 7
$ws.Range("B7").Value = 42
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 36.9
$ws.Range("E7").Value = 35.7
$ws.Range("F7").Value = 35.6
$ws.Range("G7").Value = 36.7
$ws.Range("H7").Value = 38.8
$ws.Range("I7").Value = 42.1
$ws.Range("J7").Value = 46.3
$ws.Range("K7").Value = 51.2
$ws.Range("L7").Value = 56.5

# Row 8
$ws.Range("B8").Value = 34
$ws.Range("C8").Value = 31.2
$ws.Range("D8").Value = 29.6
$ws.Range("E8").Value = 29.2
$ws.Range("F8").Value = 30.1
$ws.Range("G8").Value = 32.3
$ws.Range("H8").Value = 35.8
$ws.Range("I8").Value = 40.4
$ws.Range("J8").Value = 45.8
$ws.Range("K8").Value = 51.8
$ws.Range("L8").Value = 57.9

# Row 9
$ws.Range("B9").Value = 33.2
$ws.Range("C9").Value = 33
$ws.Range("D9").Value = 33.3
$ws.Range("E9").Value = 34.1
$ws.Range("F9").Value = 35.2
$ws.Range("G9").Value = 36.7
$ws.Range("H9").Value = 38.5
$ws.Range("I9").Value = 40.7
$ws.Range("J9").Value = 43.2
$ws.Range("K9").Value = 46.1
$ws.Range("L9").Value = 49.3

# Row 10
$ws.Range("B10").Value = 47.9
$ws.Range("C10").Value = 48.3
$ws.Range("D10").Value = 49.1
$ws.Range("E10").Value = 50.2
$ws.Range("F10").Value = 51.4
$ws.Range("G10").Value = 52.7
$ws.Range("H10").Value = 54.1
$ws.Range("I10").Value = 55.3
$ws.Range("J10").Value = 56.5
$ws.Range("K10").Value = 57.6
$ws.Range("L10").Value = 58.7

# Row 11
$ws.Range("B11").Value = 38.2
$ws.Range("C11").Value = 38.6
$ws.Range("D11").Value = 39.8
$ws.Range("E11").Value = 41.4
$ws.Range("F11").Value = 43.1
$ws.Range("G11").Value = 44.7
$ws.Range("H11").Value = 46
$ws.Range("I11").Value = 47
$ws.Range("J11").Value = 47.6
$ws.Range("K11").Value = 48
$ws.Range("L11").Value = 48.4

# Row 12
$ws.Range("B12").Value = 50.9
$ws.Range("C12").Value = 49.8
$ws.Range("D12").Value = 48.8
$ws.Range("E12").Value = 48
$ws.Range("F12").Value = 47.5
$ws.Range("G12").Value = 47.4
$ws.Range("H12").Value = 47.6
$ws.Range("I12").Value = 48.3
$ws.Range("J12").Value = 49.4
$ws.Range("K12").Value = 50.7
$ws.Range("L12").Value = 52.3

